# The document currently has a single paragraph that contains one empty
# run (only rPr/rtl, no text). The target adds a brand-new run *before*
# that empty run, carrying its own formatting (dark-gray color, 10pt
# size, white highlight) and the Russian sentence
# "Управление разработкой программных средств."
#
# Directly inserting text into the existing empty run would just grow
# that run (and its formatting) instead of creating a second, distinctly
# formatted run, so instead we:
#   1. Split the paragraph in two with InsertParagraphAfter (the new,
#      leading paragraph is an exact clone of the original empty one).
#   2. Put the new sentence + formatting into that new, first paragraph.
#   3. Delete the paragraph mark that separates the two paragraphs again,
#      merging them back into a single paragraph. Because the two runs
#      have different run formatting, Word keeps them as two separate
#      <w:r> elements instead of coalescing them - giving exactly the
#      "new run, then the original empty run" structure the diff wants.

$d = $word.ActiveDocument

$start = $d.Range(0, 0)
$start.InsertParagraphAfter()

$newPara = $d.Paragraphs(1).Range
$newPara.Text = "Управление разработкой программных средств."

$newPara.Font.Color = 2236962          # RGB(0x22,0x22,0x22) -> 222222
$newPara.Font.Size = 10                 # half-points 20 -> sz/szCs
$newPara.Font.SizeBi = 10
$newPara.HighlightColorIndex = 8        # wdWhite -> <w:highlight w:val="white"/>

$paraEnd = $d.Paragraphs(1).Range.End
$mark = $d.Range($paraEnd - 1, $paraEnd)
$mark.Delete()
